$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - PARSTEI LX Equity
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.004410933559855143
$ws.Range("D2").Value = 0.0044109444566214

# Row 3 - LEF1TREU Index
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.004352063267490375
$ws.Range("D3").Value = 0.004352065352112043

# Row 4 - SX5R Index
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.3011493610815403
$ws.Range("D4").Value = 0.3011490759015035

# Row 5 - SXUSR Index
$ws.Range("B5").Value = 0.7
$ws.Range("C5").Value = 0.678644535406819
$ws.Range("D5").Value = 0.6786447681426276

# Row 6 - BEGCGA Index
$ws.Range("B6").Value = 0.3
$ws.Range("C6").Value = 0.002944360464062921
$ws.Range("D6").Value = 0.00294437259506601

# Row 7 - LEC4TREU Index
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0.004242739850924624
$ws.Range("D7").Value = 0.00424274542942027

# Row 8 - LEATTREU Index
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.004256006369307748
$ws.Range("D8").Value = 0.004256028122649053
